$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format Price and Volume columns as Text to avoid numeric auto-conversion
$ws.Range("D2:D51").NumberFormat = "@"
$ws.Range("E2:E51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '30.634.56'
$ws.Range("E2").Value = '  +1.06%  '

# Row 3
$ws.Range("D3").Value = '1.877.85'
$ws.Range("E3").Value = '  -0.09%  '

# Row 4
$ws.Range("D4").Value = '0.9991'
$ws.Range("E4").Value = '  -0.11%  '

# Row 5
$ws.Range("D5").Value = '239.07'
$ws.Range("E5").Value = '  +0.66%  '

# Row 6
$ws.Range("D6").Value = '0.9995'
$ws.Range("E6").Value = '  -0.05%  '

# Row 7
$ws.Range("D7").Value = '0.4799'
$ws.Range("E7").Value = '  -0.51%  '

# Row 8
$ws.Range("D8").Value = '0.2832'
$ws.Range("E8").Value = '  -2.00%  '

# Row 9
$ws.Range("D9").Value = '0.06533'
$ws.Range("E9").Value = '  -0.82%  '

# Row 10
$ws.Range("D10").Value = '1.961.73'
$ws.Range("E10").Value = '  +4.41%  '

# Row 11
$ws.Range("D11").Value = '0.07465'
$ws.Range("E11").Value = '  +0.92%  '

# Row 12
$ws.Range("D12").Value = '16.64'
$ws.Range("E12").Value = '  -1.71%  '

# Row 13
$ws.Range("E13").Value = '  -1.61%  '

# Row 14
$ws.Range("D14").Value = '88.90'
$ws.Range("E14").Value = '  +1.05%  '

# Row 15
$ws.Range("D15").Value = '0.6654'
$ws.Range("E15").Value = '  +0.92%  '

# Row 16
$ws.Range("D16").Value = '30.582.22'
$ws.Range("E16").Value = '  +0.99%  '

# Row 17
$ws.Range("E17").Value = '  -1.85%  '

# Row 18
$ws.Range("D18").Value = '0.9994'

# Row 19
$ws.Range("D19").Value = '0.000007613'
$ws.Range("E19").Value = '  -1.47%  '

# Row 20
$ws.Range("B20").Value = 'BitcoinCash'
$ws.Range("C20").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D20").Value = '233.25'
$ws.Range("E20").Value = '  +19.01%  '

# Row 21
$ws.Range("B21").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C21").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D21").Value = '2.110.48'
$ws.Range("E21").Value = '  -1.22%  '

# Row 22
$ws.Range("D22").Value = '5.306'
$ws.Range("E22").Value = '  -3.00%  '

# Row 23
$ws.Range("D23").Value = '0.9995'
$ws.Range("E23").Value = '  -0.03%  '

# Row 24
$ws.Range("B24").Value = 'BitDAO'
$ws.Range("C24").Value = 'https://coinranking.com/coin/N2IgQ9Xme+bitdao-bit'
$ws.Range("D24").Value = '0.3948'
$ws.Range("E24").Value = '  -6.63%  '

# Row 25
$ws.Range("B25").Value = 'Chainlink'
$ws.Range("C25").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D25").Value = '6.216'
$ws.Range("E25").Value = '  +1.05%  '

# Row 26
$ws.Range("B26").Value = 'Cosmos'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D26").Value = '9.318'
$ws.Range("E26").Value = '  -1.11%  '

# Row 27
$ws.Range("B27").Value = 'Monero'
$ws.Range("C27").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D27").Value = '167.37'
$ws.Range("E27").Value = '  +2.46%  '

# Row 28
$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").Value = '18.81'
$ws.Range("E28").Value = '  +3.23%  '

# Row 29
$ws.Range("B29").Value = 'LidoDAOToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D29").Value = '1.956'
$ws.Range("E29").Value = '  +1.46%  '

# Row 30
$ws.Range("B30").Value = 'Toncoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D30").Value = '1.453'
$ws.Range("E30").Value = '  +0.94%  '

# Row 31
$ws.Range("B31").Value = 'Stellar'
$ws.Range("C31").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D31").Value = '0.09564'
$ws.Range("E31").Value = '  +4.60%  '

# Row 32
$ws.Range("B32").Value = 'InternetComputer(DFINITY)'
$ws.Range("C32").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D32").Value = '4.317'
$ws.Range("E32").Value = '  +1.08%  '

# Row 33
$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").Value = '4.036'
$ws.Range("E33").Value = '  -0.16%  '

# Row 34
$ws.Range("B34").Value = 'Hedera'
$ws.Range("C34").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D34").Value = '0.05031'
$ws.Range("E34").Value = '  -0.35%  '

# Row 35
$ws.Range("B35").Value = 'ARBITRUM'
$ws.Range("C35").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D35").Value = '1.214'
$ws.Range("E35").Value = '  +6.65%  '

# Row 36
$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D36").Value = '0.7502'
$ws.Range("E36").Value = '  +1.26%  '

# Row 37
$ws.Range("B37").Value = 'HuobiToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D37").Value = '2.713'
$ws.Range("E37").Value = '  +0.19%  '

# Row 38
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").Value = '0.01852'
$ws.Range("E38").Value = '  +0.84%  '

# Row 39
$ws.Range("B39").Value = 'MXToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D39").Value = '2.622'
$ws.Range("E39").Value = '  -0.38%  '

# Row 40
$ws.Range("D40").Value = '0.9148'
$ws.Range("E40").Value = '  +0.05%  '

# Row 41
$ws.Range("B41").Value = 'RenderToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D41").Value = '2.079'
$ws.Range("E41").Value = '  +0.35%  '

# Row 42
$ws.Range("B42").Value = 'Quant'
$ws.Range("C42").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D42").Value = '105.89'
$ws.Range("E42").Value = '  -0.45%  '

# Row 43
$ws.Range("B43").Value = 'TheSandbox'
$ws.Range("C43").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D43").Value = '0.4280'
$ws.Range("E43").Value = '  -0.89%  '

# Row 44
$ws.Range("B44").Value = 'FraxShare'
$ws.Range("C44").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D44").Value = '5.806'
$ws.Range("E44").Value = '  -1.19%  '

# Row 45
$ws.Range("B45").Value = 'PaxDollar'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D45").Value = '1.005'
$ws.Range("E45").Value = '  +0.55%  '

# Row 46
$ws.Range("B46").Value = 'Aptos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D46").Value = '7.488'
$ws.Range("E46").Value = '  -1.79%  '

# Row 47
$ws.Range("B47").Value = 'Aave'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D47").Value = '64.58'
$ws.Range("E47").Value = '  -0.61%  '

# Row 48
$ws.Range("B48").Value = 'Algorand'
$ws.Range("C48").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D48").Value = '0.1286'
$ws.Range("E48").Value = '  -4.47%  '

# Row 49
$ws.Range("B49").Value = 'NEARProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D49").Value = '1.485'
$ws.Range("E49").Value = '  -5.14%  '

# Row 50
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").Value = '8.945'
$ws.Range("E50").Value = '  +1.13%  '

# Row 51
$ws.Range("B51").Value = 'Elrond'
$ws.Range("C51").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D51").Value = '33.80'
$ws.Range("E51").Value = '  -0.99%  '
